$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.301593443182684
$ws.Range("D2").Value = 0.20652091539819

$ws.Range("C3").Value = -0.1423096387251648
$ws.Range("D3").Value = 0.8881306875568387

$ws.Range("C4").Value = -0.2775536872456714
$ws.Range("D4").Value = 0.7839464926165656

$ws.Range("C5").Value = 1.866660502619637
$ws.Range("D5").Value = 0.07533441735242263

$ws.Range("C6").Value = -1.244212745704023
$ws.Range("D6").Value = 0.2265147570013888

$ws.Range("C7").Value = -1.031369333911698
$ws.Range("D7").Value = 0.3135734621242485

$ws.Range("C8").Value = 0.3788220551707878
$ws.Range("D8").Value = 0.7084513348688191

$ws.Range("C9").Value = -0.1063305383308633
$ws.Range("D9").Value = 0.9162837352247604

$ws.Range("C10").Value = 2.023690442181398
$ws.Range("D10").Value = 0.05531559023917265
$ws.Range("G10").Value = "No"

$ws.Range("C11").Value = 1.822998486183408
$ws.Range("D11").Value = 0.08192295687616835
